# Commit: "Change in unit of AIC and run of new results"
#
# Across every yearly worksheet ("2000".."2100") the Abiotic-resource
# Impact Category (AIC) figures stored in columns D-G of rows 5, 7 and 8
# are re-expressed in a unit that is 1,000,000x larger than before, i.e.
# every currently-non-zero value in those cells is multiplied by 1E-6
# (divided by one million). Cells that are already 0 stay 0.
#
#   Row 5 -> column E
#   Row 7 -> columns D, E, G
#   Row 8 -> column F

$wb = $excel.ActiveWorkbook

$unitFactor = 0.000001

# Map of row number -> the column letters that hold AIC figures on that row.
$targets = @{
    5 = @("E")
    7 = @("D", "E", "G")
    8 = @("F")
}

$sheetCount = $wb.Worksheets.Count

for ($i = 1; $i -le $sheetCount; $i++) {
    $ws = $wb.Worksheets.Item($i)

    foreach ($row in $targets.Keys) {
        foreach ($col in $targets[$row]) {
            $rng = $ws.Range("$col$row")
            $v = $rng.Value()
            if ($v -ne 0) {
                $rng.Value = $v * $unitFactor
            }
        }
    }
}
